{"js": "// The commit replaces the old (Russian-dated) \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u041f\u0435\u0440\u0441\u0435\u0443\u0441 ...\" sentence\n// with the new translated sentence \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 Leo: 14. \u0438 23. \u0430\u043f\u0440\u0438\u043b\u0430, 14. \u0438 23. \u043c\u0430\u0458\u0430\"\n// in the four paragraphs of the document body that contain it. Each of those\n// paragraphs originally held several runs (plus, in the first case, a\n// _Hlk514861060 bookmark) -- all of that is collapsed into a single run with\n// the new text and default (inherited) run formatting, matching the target\n// OOXML which has a bare <w:r><w:t>...</w:t></w:r>.\n\nconst newText = \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 Leo: 14. \u0438 23. \u0430\u043f\u0440\u0438\u043b\u0430, 14. \u0438 23. \u043c\u0430\u0458\u0430\";\nconst oldText = \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u041f\u0435\u0440\u0441\u0435\u0443\u0441 \u0442\u043e\u043a\u043e\u043c 2018. \u0433\u043e\u0434\u0438\u043d\u0435 \u043f\u043e\u0441\u043c\u0430\u0442\u0440\u0430\u043c\u043e 30. \u043e\u043a\u0442\u043e\u0431\u0440\u0430 - 8. \u043d\u043e\u0432\u0435\u043c\u0431\u0440\u0430 \u0438 29. \u043d\u043e\u0432\u0435\u043c\u0431\u0440\u0430 - 8. \u0434\u0435\u0446\u0435\u043c\u0431\u0440\u0430\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text.trim();\n  if (text === oldText) {\n    // Clear every run (and the bookmark, if any) out of the paragraph, then\n    // insert a single new run with no explicit formatting overrides.\n    para.clear();\n    para.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The commit replaces the old (Russian-dated) \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u041f\u0435\u0440\u0441\u0435\u0443\u0441 ...\" sentence\n# with the new translated sentence \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 Leo: 14. \u0438 23. \u0430\u043f\u0440\u0438\u043b\u0430, 14. \u0438 23. \u043c\u0430\u0458\u0430\"\n# in every paragraph of the document body that contains it. Each of those\n# paragraphs originally held several runs (plus, in the first occurrence, a\n# _Hlk514861060 bookmark around part of the text) -- all of that collapses\n# into a single run with the new text and no explicit run-formatting\n# overrides, matching the target OOXML's bare <w:r><w:t>...</w:t></w:r>.\n\n$d = $word.ActiveDocument\n\n$old = \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u041f\u0435\u0440\u0441\u0435\u0443\u0441 \u0442\u043e\u043a\u043e\u043c 2018. \u0433\u043e\u0434\u0438\u043d\u0435 \u043f\u043e\u0441\u043c\u0430\u0442\u0440\u0430\u043c\u043e 30. \u043e\u043a\u0442\u043e\u0431\u0440\u0430 - 8. \u043d\u043e\u0432\u0435\u043c\u0431\u0440\u0430 \u0438 29. \u043d\u043e\u0432\u0435\u043c\u0431\u0440\u0430 - 8. \u0434\u0435\u0446\u0435\u043c\u0431\u0440\u0430\"\n$new = \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 Leo: 14. \u0438 23. \u0430\u043f\u0440\u0438\u043b\u0430, 14. \u0438 23. \u043c\u0430\u0458\u0430\"\n\n# Remove the tracking bookmark that wraps part of the old sentence, if present.\ntry {\n    $bm = $d.Bookmarks(\"_Hlk514861060\")\n    $bm.Delete()\n} catch {\n    # bookmark not present / already removed - nothing to do\n}\n\nforeach ($p in $d.Paragraphs) {\n    $rng = $p.Range\n    $text = $rng.Text.TrimEnd([char]13, [char]12).Trim()\n    if ($text -eq $old) {\n        $paraRange = $p.Range\n        # Exclude the trailing paragraph mark so we only touch the runs.\n        $paraRange.MoveEnd(1, -1) | Out-Null\n        # Drop all existing runs (and their formatting) then type the plain\n        # replacement text so the resulting run carries no rPr overrides.\n        $paraRange.Delete()\n        $paraRange.InsertAfter($new)\n    }\n}\n"}
